# Update the "samples" worksheet with the new survey-style data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row -----------------------------------------------------------
$headers = @(
    "Row 1",
    "How motivated are you to come to work every day?",
    "How much do you feel valued and recognized for your work?",
    "How would you rate the opportunities for professional development and career opportunities in the company?",
    "Do you feel you are treated fairly and equally?",
    "How would you rate the company's salary and benefits?",
    "How transparent are decision-making processes in the company?",
    "How would you rate the leadership skills in the company?",
    "How well are new employees integrated into the company?"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Give the header cells (B1:I1) vertical-centered alignment.
$ws.Range("B1:I1").VerticalAlignment = -4108

# --- Data rows -------------------------------------------------------------
$names = @("Jack", "Marie", "Martin", "Tom", "Vanessa", "Ylvi", "Dan")
$scores = @(
    @(0, 6, 4, 7, 4, 3, 6, 7),
    @(1, 7, 5, 8, 5, 4, 5, 8),
    @(2, 5, 4, 9, 6, 2, 3, 9),
    @(3, 6, 3, 8, 4, 1, 5, 9),
    @(4, 6, 6, 9, 5, 6, 7, 8),
    @(5, 8, 8, 8, 4, 5, 7, 9),
    @(6, 9, 5, 7, 4, 4, 6, 7)
)

for ($r = 0; $r -lt $names.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 1).Value = $scores[$r][0]
    $ws.Cells.Item($row, 2).Value = $names[$r]
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($row, $c + 3).Value = $scores[$r][$c + 1]
    }
}

# --- Column widths (B:I) ---------------------------------------------------
$ws.Range("B1:I8").ColumnWidth = 12.6640625

# --- Selection --------------------------------------------------------------
$ws.Range("F17").Select()
